# Apply the crypto list refresh described in the commit diff (GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "79.007.19"
$ws.Range("E2").Value = "  +3.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.183.79"
$ws.Range("E3").Value = "  +5.13%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.86"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "632.43"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.229"
$ws.Range("E8").Value = "  +8.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").Value = "  +5.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.179.73"
$ws.Range("E10").Value = "  +5.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  +32.53%  "

$ws.Range("E12").Value = "  +2.84%  "

$ws.Range("E13").Value = "  +4.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.769.23"
$ws.Range("E14").Value = "  +5.13%  "

$ws.Range("E15").Value = "  +15.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.43"
$ws.Range("E16").Value = "  +6.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "78.911.80"
$ws.Range("E17").Value = "  +3.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.176.60"
$ws.Range("E18").Value = "  +4.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.46"
$ws.Range("E19").Value = "  +6.91%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.24"
$ws.Range("E20").Value = "  +2.44%  "

$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.97"
$ws.Range("E21").Value = "  +28.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.81"
$ws.Range("E22").Value = "  +12.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.97"
$ws.Range("E23").Value = "  +13.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.86"
$ws.Range("E24").Value = "  +6.11%  "

$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.79"
$ws.Range("E25").Value = "  +8.53%  "

$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.14"
$ws.Range("E26").Value = "  +11.51%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "76.46"
$ws.Range("E27").Value = "  +4.06%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000114"
$ws.Range("E29").Value = "  +2.92%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.94"
$ws.Range("E31").Value = "  +6.72%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.47"
$ws.Range("E32").Value = "  +2.94%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "520.35"
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.99"
$ws.Range("E34").Value = "  +1.90%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.96"
$ws.Range("E35").Value = "  +10.38%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.135"
$ws.Range("E36").Value = "  +19.32%  "

$ws.Range("B37").Value = "Cronos"
$ws.Range("C37").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.124"
$ws.Range("E37").Value = "  +16.72%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.402"
$ws.Range("E39").Value = "  +4.19%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "164.29"
$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.01"
$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "191.75"
$ws.Range("E43").Value = "  +1.70%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.42"
$ws.Range("E44").Value = "  +4.46%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.816"
$ws.Range("E45").Value = "  +13.93%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.78"
$ws.Range("E46").Value = "  +6.05%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.31"
$ws.Range("E47").Value = "  +4.08%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.43"
$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.50"
$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.29"
$ws.Range("E50").Value = "  +11.70%  "

$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.625"
$ws.Range("E51").Value = "  +2.83%  "

